# Updates the "cryptos" price/volume table with refreshed values.
# Cells in column D that look like plain numbers are written with a
# leading apostrophe (quote-prefix) to keep them as text - matching the
# original workbook where all Price values are stored as strings - and
# then their style is reset to "Normal" so no stray numeric style is
# left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.997.15"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").Value = "2.416.12"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'573.67"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "'146.10"
$ws.Range("E6").Value = "  +6.03%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").Value = "2.447.99"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("E10").Value = "  +5.62%  "
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "'5.24"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("D14").Value = "'27.32"
$ws.Range("E14").Value = "  +6.43%  "
$ws.Range("E15").Value = "  +7.97%  "
$ws.Range("D16").Value = "2.856.27"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "62.734.19"
$ws.Range("E17").Value = "  +4.96%  "
$ws.Range("D18").Value = "2.443.16"
$ws.Range("E18").Value = "  +2.85%  "
$ws.Range("D19").Value = "'7.92"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").Value = "'10.98"
$ws.Range("E20").Value = "  +4.54%  "
$ws.Range("D21").Value = "'328.79"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("E23").Value = "  +13.27%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'65.58"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").Value = "'635.27"
$ws.Range("E26").Value = "  +14.20%  "
$ws.Range("D27").Value = "'8.51"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").Value = "0.0₃0991"
$ws.Range("E28").Value = "  +7.35%  "
$ws.Range("D30").Value = "'8.22"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").Value = "  +9.08%  "
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +4.84%  "
$ws.Range("E37").Value = "  +2.44%  "
$ws.Range("D38").Value = "'152.98"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").Value = "'5.41"
$ws.Range("E39").Value = "  +8.42%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'18.69"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  +15.49%  "
$ws.Range("E42").Value = "  +8.01%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "0.0₆0289"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").Value = "'145.01"
$ws.Range("E45").Value = "  +3.81%  "
$ws.Range("D46").Value = "'3.60"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").Value = "'20.50"
$ws.Range("E47").Value = "  +7.74%  "
$ws.Range("D48").Value = "'0.602"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").Value = "'13.04"
$ws.Range("E50").Value = "  +11.72%  "
$ws.Range("D51").Value = "'0.0921"
$ws.Range("E51").Value = "  +2.47%  "

# Reset style on cells that were forced to text via quote-prefix,
# so no visible style index is added to those cells.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
